$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.736.85"
$ws.Range("D3").Value = "3.782.74"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +1.44%  "
$ws.Range("D5").Value = "599.16"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "162.89"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("D7").Value = "3.782.26"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "6.60"
$ws.Range("E12").Value = "  +4.74%  "
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").Value = "35.03"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "4.416.57"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "3.750.97"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "67.775.51"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "18.15"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "7.00"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "458.14"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "9.46"
$ws.Range("E22").Value = "  -4.46%  "
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").Value = "82.88"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "0.0000141"
$ws.Range("E25").Value = "  -6.22%  "
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").Value = "3.927.18"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "2.20"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "7.22"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  -6.65%  "
$ws.Range("D34").Value = "28.83"
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "8.94"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "0.0992"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  +4.53%  "
$ws.Range("D39").Value = "5.78"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "0.976"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  -6.49%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D44").Value = "43.49"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D45").Value = "47.20"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").Value = "151.17"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").Value = "8.27"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "1.35"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "1.83"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "385.90"
$ws.Range("E51").Value = "  -1.79%  "
